$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and 1h volume-change data to latest scraped values

$ws.Range("D2").Value = "52.464.75"
$ws.Range("E2").Value = "  +1.02%  "

$ws.Range("D3").Value = "3.021.50"
$ws.Range("E3").Value = "  +2.16%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "356.07"
$ws.Range("E5").Value = "  +0.74%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.95"
$ws.Range("E6").Value = "  -2.73%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.565"
$ws.Range("E7").Value = "  +0.51%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.625"
$ws.Range("E9").Value = "  -1.21%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.68"
$ws.Range("E10").Value = "  -2.60%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.139"
$ws.Range("E11").Value = "  +2.09%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0863"
$ws.Range("E12").Value = "  -3.93%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.42"
$ws.Range("E13").Value = "  -2.41%  "

$ws.Range("D14").Value = "3.493.89"
$ws.Range("E14").Value = "  +1.98%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.75"
$ws.Range("E15").Value = "  -4.39%  "

$ws.Range("D16").Value = "3.013.99"
$ws.Range("E16").Value = "  +2.36%  "

$ws.Range("E17").Value = "  +2.77%  "

$ws.Range("D18").Value = "52.536.77"
$ws.Range("E18").Value = "  +0.84%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.55"
$ws.Range("E19").Value = "  +8.72%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.57"
$ws.Range("E20").Value = "  -2.35%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.77"
$ws.Range("E21").Value = "  -4.63%  "

$ws.Range("D22").Value = "0.0₃0977"
$ws.Range("E22").Value = "  -1.64%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.82"
$ws.Range("E23").Value = "  -2.53%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "265.96"
$ws.Range("E24").Value = "  -2.32%  "

$ws.Range("E25").Value = "  -1.99%  "

$ws.Range("E26").Value = "  -1.26%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.21"
$ws.Range("E27").Value = "  -1.16%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.67"
$ws.Range("E28").Value = "  +2.56%  "

$ws.Range("E29").Value = "  +0.00%  "

$ws.Range("E30").Value = "  -1.89%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.43"
$ws.Range("E31").Value = "  +0.80%  "

$ws.Range("E32").Value = "  -3.07%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "36.86"
$ws.Range("E33").Value = "  -2.55%  "

$ws.Range("E34").Value = "  +17.00%  "

$ws.Range("E35").Value = "  -3.98%  "

$ws.Range("E36").Value = "  -1.19%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  +0.02%  "

$ws.Range("E38").Value = "  -4.62%  "

$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.02"
$ws.Range("E39").Value = "  -2.01%  "

$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.10"
$ws.Range("E40").Value = "  -4.39%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.72"
$ws.Range("E41").Value = "  +1.31%  "

$ws.Range("E42").Value = "  -0.71%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.08"
$ws.Range("E43").Value = "  -3.27%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "123.57"
$ws.Range("E44").Value = "  +8.13%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.16"
$ws.Range("E45").Value = "  -0.92%  "

$ws.Range("D46").Value = "2.138.34"
$ws.Range("E46").Value = "  -1.59%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.42"
$ws.Range("E47").Value = "  -4.35%  "

$ws.Range("E48").Value = "  -5.24%  "

$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "3.311.21"
$ws.Range("E49").Value = "  +1.79%  "

$ws.Range("B50").Value = "TheGraph"
$ws.Range("C50").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.248"
$ws.Range("E50").Value = "  +1.22%  "

$ws.Range("B51").Value = "BEAM"
$ws.Range("C51").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0337"
$ws.Range("E51").Value = "  -1.14%  "

